$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '68.222.93'
$ws.Range('E2').Value = '  -0.41%  '

# Row 3
$ws.Range('D3').Value = '3.887.98'
$ws.Range('E3').Value = '  -0.80%  '

# Row 4
$ws.Range('E4').Value = '  -0.01%  '

# Row 5
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '482.63'
$ws.Range('E5').Value = '  -0.47%  '

# Row 6
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '145.41'
$ws.Range('E6').Value = '  -1.27%  '

# Row 7
$ws.Range('E7').Value = '  +0.39%  '

# Row 8
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.997'
$ws.Range('E8').Value = '  -0.12%  '

# Row 9
$ws.Range('E9').Value = '  +2.70%  '

# Row 10
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.182'
$ws.Range('E10').Value = '  +8.44%  '

# Row 11
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0000354'
$ws.Range('E11').Value = '  +0.09%  '

# Row 12
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '43.02'
$ws.Range('E12').Value = '  +1.43%  '

# Row 13
$ws.Range('E13').Value = '  +0.00%  '

# Row 14
$ws.Range('D14').Value = '4.506.08'
$ws.Range('E14').Value = '  -0.90%  '

# Row 15
$ws.Range('D15').Value = '3.900.72'
$ws.Range('E15').Value = '  -0.26%  '

# Row 16
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '14.22'
$ws.Range('E16').Value = '  -2.42%  '

# Row 17
$ws.Range('E17').Value = '  -0.67%  '

# Row 18
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '19.91'
$ws.Range('E18').Value = '  +1.02%  '

# Row 19
$ws.Range('E19').Value = '  +0.30%  '

# Row 20
$ws.Range('D20').Value = '68.252.05'
$ws.Range('E20').Value = '  -0.65%  '

# Row 21
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '429.07'
$ws.Range('E21').Value = '  -0.77%  '

# Row 22
$ws.Range('E22').Value = '  +8.65%  '

# Row 23
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '14.79'
$ws.Range('E23').Value = '  +1.84%  '

# Row 24
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '12.41'
$ws.Range('E24').Value = '  +15.72%  '

# Row 25
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '88.67'
$ws.Range('E25').Value = '  +2.09%  '

# Row 26
$ws.Range('E26').Value = '  +2.25%  '

# Row 27
$ws.Range('E27').Value = '  -3.55%  '

# Row 28
$ws.Range('E28').Value = '  -1.94%  '

# Row 29
$ws.Range('E29').Value = '  -3.43%  '

# Row 30
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '717.52'
$ws.Range('E30').Value = '  -0.08%  '

# Row 31
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '13.49'
$ws.Range('E31').Value = '  +2.03%  '

# Row 32
$ws.Range('E32').Value = '  +0.34%  '

# Row 33
$ws.Range('E33').Value = '  +2.79%  '

# Row 34
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '61.91'
$ws.Range('E34').Value = '  +6.06%  '

# Row 35
$ws.Range('B35').Value = 'NEARProtocol'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '6.06'
$ws.Range('E35').Value = '  +9.63%  '

# Row 36
$ws.Range('B36').Value = 'PEPE'
$ws.Range('C36').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D36').Value = '0.0₃0873'
$ws.Range('E36').Value = '  -3.32%  '

# Row 37
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '40.87'
$ws.Range('E37').Value = '  -1.21%  '

# Row 38
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.400'
$ws.Range('E38').Value = '  +16.80%  '

# Row 39
$ws.Range('E39').Value = '  -3.26%  '

# Row 40
$ws.Range('B40').Value = 'Dai'
$ws.Range('C40').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.997'
$ws.Range('E40').Value = '  -0.17%  '

# Row 41
$ws.Range('B41').Value = 'VeChain'
$ws.Range('C41').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.0498'
$ws.Range('E41').Value = '  +6.14%  '

# Row 42
$ws.Range('B42').Value = 'Fetch.AI'
$ws.Range('C42').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '3.01'
$ws.Range('E42').Value = '  +5.63%  '

# Row 43
$ws.Range('E43').Value = '  +3.18%  '

# Row 44
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '2.96'
$ws.Range('E44').Value = '  -3.19%  '

# Row 45
$ws.Range('E45').Value = '  +2.04%  '

# Row 46
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '3.36'
$ws.Range('E46').Value = '  +5.51%  '

# Row 47
$ws.Range('E47').Value = '  -0.06%  '

# Row 48
$ws.Range('B48').Value = 'BabyDogeCoin'
$ws.Range('C48').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D48').Value = '0.0₆0354'
$ws.Range('E48').Value = '  +29.14%  '

# Row 49
$ws.Range('B49').Value = 'LidoDAOToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '3.37'
$ws.Range('E49').Value = '  -0.92%  '

# Row 50
$ws.Range('E50').Value = '  -2.31%  '

# Row 51
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '144.48'
$ws.Range('E51').Value = '  -2.38%  '
